$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update search result data: B5 changes from "Butter" to "Milk"
$ws.Range("B5").Value = "Milk"

# Update the active cell selection to C5
$ws.Range("C5").Select()
